$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# Slide 16: the content placeholder's last paragraph was an empty
# bullet (just an endParaRPr). Fill it with "Pruebas funcionales"
# (as three runs: "Pruebas" / " " / "funcionales", sz=3200) and add
# a new empty bullet paragraph after it so the placeholder still ends
# with a blank bullet line, same as before the edit.
# ------------------------------------------------------------------
$s16 = $p.Slides.Item(16)
$shp16 = $s16.Shapes.Item(2)
$tr16 = $shp16.TextFrame.TextRange

$para3 = $tr16.Paragraphs(3, 1)
$para3.Text = "Pruebas funcionales"

# Create a new trailing empty bullet paragraph (re-using a throwaway
# placeholder character so we can cleanly delete it afterwards and
# leave the new paragraph with no run at all, matching the original
# trailing-empty-paragraph shape).
$para3.InsertAfter("`rX")
$para4 = $tr16.Paragraphs(4, 1)
$placeholder = $para4.Characters(1, 1)
$placeholder.Delete()

# Split "Pruebas funcionales" into its three runs and size them to
# match the surrounding bullet text (32pt).
$para3 = $tr16.Paragraphs(3, 1)
$run1 = $para3.Characters(1, 7)
$run1.Font.Size = 32
$run2 = $para3.Characters(8, 1)
$run2.Font.Size = 32
$run3 = $para3.Characters(9, 11)
$run3.Font.Size = 32

# ------------------------------------------------------------------
# Slide 2: "Pantall" + "a" were two separate runs that should be a
# single "Pantalla" run (typo fix / run merge).
# ------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shp2 = $s2.Shapes.Item(2)
$tr2 = $shp2.TextFrame.TextRange

$para = $tr2.Paragraphs(4, 1)
$word = $para.Characters(1, 8)
$word.Text = "Pantalla"
